$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1140
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -7340

$ws.Range("H137").Value = 1897.5
$ws.Range("I137").Value = 1726
$ws.Range("K137").Value = 5178
$ws.Range("M137").Value = -2628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1987.4783
$ws.Range("I61").Value = 1403.2142
$ws.Range("K61").Value = 1403.2142
$ws.Range("M61").Value = -1191.2142

$ws.Range("H110").Value = 1969.3334
$ws.Range("I110").Value = 2028
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 2028
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 17
$ws.Range("N110").Value = -5590

$ws.Range("H122").Value = 1345.3572
$ws.Range("I122").Value = 1235.421
$ws.Range("J122").Value = 1577.4445
$ws.Range("K122").Value = 3706.263
$ws.Range("L122").Value = 4732.333500000001
$ws.Range("M122").Value = -1256.263
$ws.Range("N122").Value = -9632.333500000001

$ws.Range("H136").Value = 1987.4783
$ws.Range("I136").Value = 1403.2142
$ws.Range("K136").Value = 4209.642599999999
$ws.Range("M136").Value = -1659.642599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 35000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 35000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 35000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -36372

$ws.Range("H66").Value = 35000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 35000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 105000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -111864

$ws.Range("H99").Value = 1791.2142
$ws.Range("I99").Value = 1606.875
$ws.Range("J99").Value = 2037
$ws.Range("K99").Value = 1606.875
$ws.Range("L99").Value = 2037
$ws.Range("M99").Value = -108.875
$ws.Range("N99").Value = -5033

$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws.Range("H134").Value = 2236
$ws.Range("I134").Value = 2191.6
$ws.Range("J134").Value = 2288.2354
$ws.Range("K134").Value = 6574.799999999999
$ws.Range("L134").Value = 6864.706200000001
$ws.Range("M134").Value = -4039.799999999999
$ws.Range("N134").Value = -11934.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 15000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16996

$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 15000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -54984

$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080

$ws.Range("H122").Value = 1311.7
$ws.Range("J122").Value = 1938
$ws.Range("L122").Value = 5814
$ws.Range("N122").Value = -10714

$ws.Range("H132").Value = 3292188.5
$ws.Range("I132").Value = 2202.7
$ws.Range("J132").Value = 6947728.5
$ws.Range("K132").Value = 6608.099999999999
$ws.Range("L132").Value = 20843185.5
$ws.Range("M132").Value = -4078.099999999999
$ws.Range("N132").Value = -20848245.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 30.2
$ws.Range("I14").Value = 30.2
$ws.Range("K14").Value = 90.59999999999999
$ws.Range("M14").Value = 82.40000000000001

$ws.Range("I92").Value = 200
$ws.Range("K92").Value = 600
$ws.Range("M92").Value = 648

$ws.Range("H107").Value = 9540.299999999999
$ws.Range("I107").Value = 1100
$ws.Range("J107").Value = 11650.375
$ws.Range("K107").Value = 3300
$ws.Range("L107").Value = 34951.125
$ws.Range("M107").Value = -1380
$ws.Range("N107").Value = -38791.125

$ws.Range("H131").Value = 768.1828
$ws.Range("J131").Value = 801.8941
$ws.Range("L131").Value = 2405.6823
$ws.Range("N131").Value = -12485.6823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 17857956
$ws.Range("I113").Value = 31250648
$ws.Range("J113").Value = 1032.1666
$ws.Range("K113").Value = 31250648
$ws.Range("L113").Value = 1032.1666
$ws.Range("M113").Value = -31248478
$ws.Range("N113").Value = -5372.1666

$ws.Range("H126").Value = 3740
$ws.Range("I126").Value = 5300
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 15900
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -13430
$ws.Range("N126").Value = -9140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 40856.668
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 40856.668
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 40856.668
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -46846.668

$ws.Range("H132").Value = 8088.4287
$ws.Range("I132").Value = 9646.619000000001
$ws.Range("J132").Value = 3413.8572
$ws.Range("K132").Value = 28939.857
$ws.Range("L132").Value = 10241.5716
$ws.Range("M132").Value = -26409.857
$ws.Range("N132").Value = -15301.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 20813
$ws.Range("J27").Value = 20813
$ws.Range("L27").Value = 20813
$ws.Range("N27").Value = -20951

$ws.Range("H107").Value = 45833824
$ws.Range("I107").Value = 41667068
$ws.Range("J107").Value = 47619576
$ws.Range("K107").Value = 125001204
$ws.Range("L107").Value = 142858728
$ws.Range("M107").Value = -124999284
$ws.Range("N107").Value = -142862568

$ws.Range("H115").Value = 40799
$ws.Range("J115").Value = 40799
$ws.Range("L115").Value = 40799
$ws.Range("N115").Value = -43933

$ws.Range("H122").Value = 2023.8125
$ws.Range("I122").Value = 1675.4286
$ws.Range("K122").Value = 5026.2858
$ws.Range("M122").Value = -2576.2858

$ws.Range("H136").Value = 1699.6177
$ws.Range("I136").Value = 1515.742
$ws.Range("J136").Value = 3599.6667
$ws.Range("K136").Value = 4547.226
$ws.Range("L136").Value = 10799.0001
$ws.Range("M136").Value = -1997.226
$ws.Range("N136").Value = -15899.0001
